$d = $word.ActiveDocument

# Helper: strip the <w:numPr>...</w:numPr> block from a paragraph's
# formatting by round-tripping the paragraph's OOXML through
# WordOpenXML / InsertXML (there is no direct "delete this element" call
# on ListFormat - RemoveNumbers() only zeroes w:numId).
function Remove-NumPr($para) {
    $r = $para.Range
    $full = $r.WordOpenXML
    $modified = $full -replace '<w:numPr>.*?</w:numPr>', ''
    $r.InsertXML($modified)
}

# 1) Drop the "Descripción: " lead-in from the three numbered list items.
#    (Replace across the whole document - Range.Text assignment would eat
#    the paragraph mark, so use Find/Replace instead.)
$d.Content.Find.Execute("Descripción: ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# 2) Remove the three short bold "<Actor>-<Historia>:" lead-in paragraphs
#    that sit just above each numbered item. Deleting paragraph 1 three
#    times in a row (re-fetching the collection each time) always lands on
#    the next remaining header because the numbered item that used to
#    follow it becomes paragraph 1 and is skipped over by index 2.
$d.Paragraphs(1).Range.Delete()
$d.Paragraphs(2).Range.Delete()
$d.Paragraphs(3).Range.Delete()

# 3) The three remaining numbered paragraphs (now paragraphs 1, 2 and 3)
#    lose their numbering entirely (the <w:numPr> element disappears, not
#    just its w:numId).
Remove-NumPr($d.Paragraphs(1))
Remove-NumPr($d.Paragraphs(2))
Remove-NumPr($d.Paragraphs(3))
